# [1N] Redimension feature added. Some minor bugs to fix regarding the Absolute resizing.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Section 1a ("1n") progress bumped to fully complete.
$ws.Range("D20").Value = 100

# Section 1o progress now partially done (was stuck at 0).
$ws.Range("D21").Value = 80

# Note left about remaining bugs in the Absolute resizing for that task.
$ws.Range("G21").Value = "alguns bugs"

# Section 1u finished.
$ws.Range("D28").Value = 100

# Widen the new "notes" column so the comment is readable.
$ws.Columns.Item(7).AutoFit()

# Leave the cursor/selection where the author was last working.
$ws.Range("H21").Select() | Out-Null
